# Update EIA Table A.7.B for the 2017-01-31 monthly refresh
# (Year-to-Date through October 2016 -> November 2016), chunk 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subtitle text: October -> November 2016
$ws.Range("A2").Value = "by End-Use Sector, Census Division, and State, Year-to-Date through November 2016"

# Updated Relative Standard Error data values
$ws.Range("C4").Value = 0.37
$ws.Range("B6").Value = 0.27
$ws.Range("B7").Value = 0.36
$ws.Range("B8").Value = 0.27
$ws.Range("B12").Value = 0.15
$ws.Range("B14").Value = 0.17
$ws.Range("C14").Value = 0.32
$ws.Range("B15").Value = 0.2
$ws.Range("B16").Value = 0.34
$ws.Range("B18").Value = 0.23
$ws.Range("B19").Value = 0.36
$ws.Range("B20").Value = 0.46
$ws.Range("B21").Value = 0.36
$ws.Range("B29").Value = 0.31
$ws.Range("C31").Value = 0.4
$ws.Range("F31").Value = 0.3
$ws.Range("B32").Value = 0.4
$ws.Range("F32").Value = 0.36
$ws.Range("B34").Value = 0.2
$ws.Range("B38").Value = 0.2
$ws.Range("D38").Value = 0.27
$ws.Range("B39").Value = 0.45
$ws.Range("B44").Value = 0.43
$ws.Range("B48").Value = 0.47
$ws.Range("B49").Value = 0.26
$ws.Range("B50").Value = 0.28000000000000003
$ws.Range("D51").Value = 5
$ws.Range("B54").Value = 0.27
$ws.Range("B58").Value = 0.15
$ws.Range("F59").Value = 0.4
$ws.Range("B61").Value = 0.47
